$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Restricciones de Usuario")
$ws.Range("A1").Value = "HELLO"
Write-Output $ws.Range("A1").Value2
